$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

$ws.Range("R1").Value = "TErm tyPe"
$ws.Range("R2").Value = "Germplasm AtTRIBUTE"

$ws.Range("R4").Select()
